$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text representation exactly (avoid numeric auto-conversion)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.856.39'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.229.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.20'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +6.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.630'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.90'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.12%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.593'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +5.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.36'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +15.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0969'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.16'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.18'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.106'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.564.80'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.99'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.864'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.229.97'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.893.78'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0967'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.22'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.97'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.10'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +7.19%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +11.58%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.53'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.72'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.80%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.40'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.73'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.126'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.125'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.56'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0734'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.74'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.13'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +20.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.94'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +8.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0303'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +12.71%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.00'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '67.32'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.01'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +18.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.01'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.203'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.80%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.09%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.65'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.50%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.16'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +7.04%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.09%  '
